$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new review row (row 7): same appid/keyword as the existing rows,
# a new reviewer email, reusing the existing "cybworking@gmail.com" recovery
# address and the "27/5/2019 15:59" timestamp from earlier rows, and a new
# review comment.
$ws.Range("A7").Value = "com.singleton.strechy"
$ws.Range("B7").Value = "stretchy"
$ws.Range("C7").Value = "sm6502345@gmail.com"
$ws.Range("D7").Value = "cybworking@gmail.com"
$ws.Range("E7").Value = "27/5/2019 15:59"
$ws.Range("F7").Value = "I think I finish the game! Please more levels! Hurry up!"

# Mailto hyperlinks on the email / recovery-email cells, same as the
# existing rows above.
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:sm6502345@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "sm6502345@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:cybworking@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "cybworking@gmail.com")

# Match the formatting used by the rest of the data rows by copying the
# formats (only) from row 6 down onto the new row — this reuses the same
# cell styles instead of minting new ones.
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F7").Select() | Out-Null
